$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.017.93'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '3.385.09'
$ws.Range("E3").Value = '  -2.06%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.42'
$ws.Range("E5").Value = '  -1.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.27'
$ws.Range("E6").Value = '  -3.04%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.606'
$ws.Range("E8").Value = '  +2.93%  '
$ws.Range("D9").Value = '3.384.98'
$ws.Range("E9").Value = '  -2.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.14'
$ws.Range("E10").Value = '  -2.93%  '
$ws.Range("E11").Value = '  -4.08%  '
$ws.Range("E12").Value = '  -1.46%  '
$ws.Range("D13").Value = '3.971.39'
$ws.Range("E13").Value = '  -1.80%  '
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("E15").Value = '  -4.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.63'
$ws.Range("E16").Value = '  -4.21%  '
$ws.Range("D17").Value = '64.111.63'
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").Value = '3.382.90'
$ws.Range("E18").Value = '  -1.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.31'
$ws.Range("E19").Value = '  -1.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.84'
$ws.Range("E20").Value = '  -4.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '372.97'
$ws.Range("E21").Value = '  -3.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.91'
$ws.Range("E22").Value = '  -3.90%  '
$ws.Range("E23").Value = '  -0.37%  '
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '71.42'
$ws.Range("E25").Value = '  -2.36%  '
$ws.Range("E26").Value = '  -4.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.83'
$ws.Range("E27").Value = '  +3.44%  '
$ws.Range("E28").Value = '  -2.68%  '
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.45'
$ws.Range("E30").Value = '  +0.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.02'
$ws.Range("E31").Value = '  -2.62%  '
$ws.Range("E32").Value = '  -1.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.01'
$ws.Range("E33").Value = '  -3.17%  '
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("E35").Value = '  +4.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '159.90'
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.87'
$ws.Range("E37").Value = '  -0.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0754'
$ws.Range("E38").Value = '  -3.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.60'
$ws.Range("E39").Value = '  -3.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.66'
$ws.Range("E40").Value = '  +0.84%  '
$ws.Range("D41").Value = '2.824.14'
$ws.Range("E41").Value = '  -3.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.55'
$ws.Range("E42").Value = '  +1.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.42'
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0312'
$ws.Range("E44").Value = '  -2.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.759'
$ws.Range("E45").Value = '  -1.41%  '
$ws.Range("E46").Value = '  +6.47%  '
$ws.Range("E47").Value = '  -2.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '309.61'
$ws.Range("E48").Value = '  +4.16%  '
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("E50").Value = '  -0.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.849'
$ws.Range("E51").Value = '  -1.11%  '
